$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 40.75339133333333
$ws.Cells.Item(2, 8).Value = 122.260174
$ws.Cells.Item(2, 9).Value = 0.02126536631186857
$ws.Cells.Item(2, 10).Value = 0.02126536631186857
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 239.0839323333333
$ws.Cells.Item(2, 14).Value = 717.251797
$ws.Cells.Item(2, 15).Value = 0.4086975387666237
$ws.Cells.Item(2, 16).Value = 0.4086975387666237
$ws.Cells.Item(2, 17).Value = 9743.48105589252
$ws.Cells.Item(2, 18).Value = 87691.32950303268
$ws.Cells.Item(2, 19).Value = 0.008691102872631357
$ws.Cells.Item(2, 20).Value = 0.008691102872631357
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 40.75339133333333
$ws.Cells.Item(3, 8).Value = 122.260174
$ws.Cells.Item(3, 9).Value = 0.02126536631186857
$ws.Cells.Item(3, 10).Value = 0.02126536631186857
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 117.0512696666667
$ws.Cells.Item(3, 14).Value = 351.153809
$ws.Cells.Item(3, 15).Value = 0.2000910950200451
$ws.Cells.Item(3, 16).Value = 0.2000910950200451
$ws.Cells.Item(3, 17).Value = 4770.236198789196
$ws.Cells.Item(3, 18).Value = 42932.12578910277
$ws.Cells.Item(3, 19).Value = 0.004255010431344159
$ws.Cells.Item(3, 20).Value = 0.004255010431344159
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 40.75339133333333
$ws.Cells.Item(4, 8).Value = 122.260174
$ws.Cells.Item(4, 9).Value = 0.02126536631186857
$ws.Cells.Item(4, 10).Value = 0.02126536631186857
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 171.15883
$ws.Cells.Item(4, 14).Value = 513.47649
$ws.Cells.Item(4, 15).Value = 0.2925842480357353
$ws.Cells.Item(4, 16).Value = 0.2925842480357353
$ws.Cells.Item(4, 17).Value = 6975.302779145473
$ws.Cells.Item(4, 18).Value = 62777.72501230927
$ws.Cells.Item(4, 19).Value = 0.006221911211562522
$ws.Cells.Item(4, 20).Value = 0.006221911211562522
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 40.75339133333333
$ws.Cells.Item(5, 8).Value = 122.260174
$ws.Cells.Item(5, 9).Value = 0.02126536631186857
$ws.Cells.Item(5, 10).Value = 0.02126536631186857
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 57.695868
$ws.Cells.Item(5, 14).Value = 173.087604
$ws.Cells.Item(5, 15).Value = 0.09862711817759588
$ws.Cells.Item(5, 16).Value = 0.09862711817759588
$ws.Cells.Item(5, 17).Value = 2351.302286920344
$ws.Cells.Item(5, 18).Value = 21161.7205822831
$ws.Cells.Item(5, 19).Value = 0.002097341796330528
$ws.Cells.Item(5, 20).Value = 0.002097341796330528
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 1689.289306666667
$ws.Cells.Item(6, 8).Value = 5067.86792
$ws.Cells.Item(6, 9).Value = 0.8814813868902838
$ws.Cells.Item(6, 10).Value = 0.8814813868902838
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 239.0839323333333
$ws.Cells.Item(6, 14).Value = 717.251797
$ws.Cells.Item(6, 15).Value = 0.4086975387666237
$ws.Cells.Item(6, 16).Value = 0.4086975387666237
$ws.Cells.Item(6, 17).Value = 403881.9302865169
$ws.Cells.Item(6, 18).Value = 3634937.372578652
$ws.Cells.Item(6, 19).Value = 0.3602592732906489
$ws.Cells.Item(6, 20).Value = 0.3602592732906489
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 1689.289306666667
$ws.Cells.Item(7, 8).Value = 5067.86792
$ws.Cells.Item(7, 9).Value = 0.8814813868902838
$ws.Cells.Item(7, 10).Value = 0.8814813868902838
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 117.0512696666667
$ws.Cells.Item(7, 14).Value = 351.153809
$ws.Cells.Item(7, 15).Value = 0.2000910950200451
$ws.Cells.Item(7, 16).Value = 0.2000910950200451
$ws.Cells.Item(7, 17).Value = 197733.4581796564
$ws.Cells.Item(7, 18).Value = 1779601.123616907
$ws.Cells.Item(7, 19).Value = 0.1763765759426649
$ws.Cells.Item(7, 20).Value = 0.1763765759426649
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 1689.289306666667
$ws.Cells.Item(8, 8).Value = 5067.86792
$ws.Cells.Item(8, 9).Value = 0.8814813868902838
$ws.Cells.Item(8, 10).Value = 0.8814813868902838
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 171.15883
$ws.Cells.Item(8, 14).Value = 513.47649
$ws.Cells.Item(8, 15).Value = 0.2925842480357353
$ws.Cells.Item(8, 16).Value = 0.2925842480357353
$ws.Cells.Item(8, 17).Value = 289136.7812605778
$ws.Cells.Item(8, 18).Value = 2602231.031345201
$ws.Cells.Item(8, 19).Value = 0.2579075687407907
$ws.Cells.Item(8, 20).Value = 0.2579075687407907
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 1689.289306666667
$ws.Cells.Item(9, 8).Value = 5067.86792
$ws.Cells.Item(9, 9).Value = 0.8814813868902838
$ws.Cells.Item(9, 10).Value = 0.8814813868902838
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 57.695868
$ws.Cells.Item(9, 14).Value = 173.087604
$ws.Cells.Item(9, 15).Value = 0.09862711817759588
$ws.Cells.Item(9, 16).Value = 0.09862711817759588
$ws.Cells.Item(9, 17).Value = 97465.01285125151
$ws.Cells.Item(9, 18).Value = 877185.1156612636
$ws.Cells.Item(9, 19).Value = 0.08693796891617914
$ws.Cells.Item(9, 20).Value = 0.08693796891617914
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 100.9654023333333
$ws.Cells.Item(10, 8).Value = 302.896207
$ws.Cells.Item(10, 9).Value = 0.05268435816499466
$ws.Cells.Item(10, 10).Value = 0.05268435816499466
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 239.0839323333333
$ws.Cells.Item(10, 14).Value = 717.251797
$ws.Cells.Item(10, 15).Value = 0.4086975387666237
$ws.Cells.Item(10, 16).Value = 0.4086975387666237
$ws.Cells.Item(10, 17).Value = 24139.20541947044
$ws.Cells.Item(10, 18).Value = 217252.848775234
$ws.Cells.Item(10, 19).Value = 0.02153196751353259
$ws.Cells.Item(10, 20).Value = 0.02153196751353259
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 100.9654023333333
$ws.Cells.Item(11, 8).Value = 302.896207
$ws.Cells.Item(11, 9).Value = 0.05268435816499466
$ws.Cells.Item(11, 10).Value = 0.05268435816499466
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 117.0512696666667
$ws.Cells.Item(11, 14).Value = 351.153809
$ws.Cells.Item(11, 15).Value = 0.2000910950200451
$ws.Cells.Item(11, 16).Value = 0.2000910950200451
$ws.Cells.Item(11, 17).Value = 11818.1285355225
$ws.Cells.Item(11, 18).Value = 106363.1568197025
$ws.Cells.Item(11, 19).Value = 0.01054167091566203
$ws.Cells.Item(11, 20).Value = 0.01054167091566203
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 100.9654023333333
$ws.Cells.Item(12, 8).Value = 302.896207
$ws.Cells.Item(12, 9).Value = 0.05268435816499466
$ws.Cells.Item(12, 10).Value = 0.05268435816499466
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 171.15883
$ws.Cells.Item(12, 14).Value = 513.47649
$ws.Cells.Item(12, 15).Value = 0.2925842480357353
$ws.Cells.Item(12, 16).Value = 0.2925842480357353
$ws.Cells.Item(12, 17).Value = 17281.1201338526
$ws.Cells.Item(12, 18).Value = 155530.0812046734
$ws.Cells.Item(12, 19).Value = 0.01541461331695031
$ws.Cells.Item(12, 20).Value = 0.01541461331695031
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 100.9654023333333
$ws.Cells.Item(13, 8).Value = 302.896207
$ws.Cells.Item(13, 9).Value = 0.05268435816499466
$ws.Cells.Item(13, 10).Value = 0.05268435816499466
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 57.695868
$ws.Cells.Item(13, 14).Value = 173.087604
$ws.Cells.Item(13, 15).Value = 0.09862711817759588
$ws.Cells.Item(13, 16).Value = 0.09862711817759588
$ws.Cells.Item(13, 17).Value = 5825.286525590892
$ws.Cells.Item(13, 18).Value = 52427.57873031803
$ws.Cells.Item(13, 19).Value = 0.005196106418849717
$ws.Cells.Item(13, 20).Value = 0.005196106418849717
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 85.41274733333334
$ws.Cells.Item(14, 8).Value = 256.238242
$ws.Cells.Item(14, 9).Value = 0.04456888863285297
$ws.Cells.Item(14, 10).Value = 0.04456888863285297
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 239.0839323333333
$ws.Cells.Item(14, 14).Value = 717.251797
$ws.Cells.Item(14, 15).Value = 0.4086975387666237
$ws.Cells.Item(14, 16).Value = 0.4086975387666237
$ws.Cells.Item(14, 17).Value = 20420.81550384677
$ws.Cells.Item(14, 18).Value = 183787.3395346209
$ws.Cells.Item(14, 19).Value = 0.01821519508981076
$ws.Cells.Item(14, 20).Value = 0.01821519508981076
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 85.41274733333334
$ws.Cells.Item(15, 8).Value = 256.238242
$ws.Cells.Item(15, 9).Value = 0.04456888863285297
$ws.Cells.Item(15, 10).Value = 0.04456888863285297
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 117.0512696666667
$ws.Cells.Item(15, 14).Value = 351.153809
$ws.Cells.Item(15, 15).Value = 0.2000910950200451
$ws.Cells.Item(15, 16).Value = 0.2000910950200451
$ws.Cells.Item(15, 17).Value = 9997.670521084865
$ws.Cells.Item(15, 18).Value = 89979.0346897638
$ws.Cells.Item(15, 19).Value = 0.00891783773037399
$ws.Cells.Item(15, 20).Value = 0.00891783773037399
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 85.41274733333334
$ws.Cells.Item(16, 8).Value = 256.238242
$ws.Cells.Item(16, 9).Value = 0.04456888863285297
$ws.Cells.Item(16, 10).Value = 0.04456888863285297
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 171.15883
$ws.Cells.Item(16, 14).Value = 513.47649
$ws.Cells.Item(16, 15).Value = 0.2925842480357353
$ws.Cells.Item(16, 16).Value = 0.2925842480357353
$ws.Cells.Item(16, 17).Value = 14619.14590065895
$ws.Cells.Item(16, 18).Value = 131572.3131059306
$ws.Cells.Item(16, 19).Value = 0.01304015476643171
$ws.Cells.Item(16, 20).Value = 0.01304015476643171
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 85.41274733333334
$ws.Cells.Item(17, 8).Value = 256.238242
$ws.Cells.Item(17, 9).Value = 0.04456888863285297
$ws.Cells.Item(17, 10).Value = 0.04456888863285297
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 57.695868
$ws.Cells.Item(17, 14).Value = 173.087604
$ws.Cells.Item(17, 15).Value = 0.09862711817759588
$ws.Cells.Item(17, 16).Value = 0.09862711817759588
$ws.Cells.Item(17, 17).Value = 4927.962595661353
$ws.Cells.Item(17, 18).Value = 44351.66336095217
$ws.Cells.Item(17, 19).Value = 0.004395701046236499
$ws.Cells.Item(17, 20).Value = 0.004395701046236499
